$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18000
$ws.Range("I21").Value = 18000
$ws.Range("K21").Value = 18000
$ws.Range("M21").Value = -17532
$ws.Range("H23").Value = 18000
$ws.Range("I23").Value = 18000
$ws.Range("K23").Value = 18000
$ws.Range("M23").Value = -17766
$ws.Range("H62").Value = 4200
$ws.Range("I62").Value = 3933.3333
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3933.3333
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3309.3333
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4200
$ws.Range("I65").Value = 3933.3333
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 19666.6665
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -16546.6665
$ws.Range("N65").Value = -31240
$ws.Range("H112").Value = 2477.5
$ws.Range("J112").Value = 2583.5
$ws.Range("L112").Value = 7750.5
$ws.Range("N112").Value = -9966.5
$ws.Range("H138").Value = 1737.2
$ws.Range("I138").Value = 750.5714
$ws.Range("J138").Value = 1999.4684
$ws.Range("K138").Value = 2251.7142
$ws.Range("L138").Value = 5998.4052
$ws.Range("M138").Value = 2888.2858
$ws.Range("N138").Value = -16278.4052
$ws.Range("H141").Value = 1534.4117
$ws.Range("I141").Value = 1567.8125
$ws.Range("K141").Value = 4703.4375
$ws.Range("M141").Value = 476.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 220.5
$ws.Range("I5").Value = 204.6
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 204.6
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -92.59999999999999
$ws.Range("N5").Value = -524
$ws.Range("H32").Value = 3546.125
$ws.Range("I32").Value = 3815.0952
$ws.Range("K32").Value = 3815.0952
$ws.Range("M32").Value = -3528.0952
$ws.Range("H74").Value = 644.6177
$ws.Range("I74").Value = 644.6177
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 644.6177
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 229.3823
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 644.6177
$ws.Range("I77").Value = 644.6177
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3223.0885
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 1144.9115
$ws.Range("N77").Value = ""
$ws.Range("H122").Value = 1006.5
$ws.Range("I122").Value = 1004
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 3012
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = -562
$ws.Range("N122").Value = -7942
$ws.Range("H139").Value = 47809.668
$ws.Range("J139").Value = 47809.668
$ws.Range("L139").Value = 47809.668
$ws.Range("N139").Value = -58089.668
$ws.Range("H140").Value = 99268
$ws.Range("J140").Value = 99268
$ws.Range("L140").Value = 99268
$ws.Range("N140").Value = -109628

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 220.5
$ws.Range("I4").Value = 204.6
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 204.6
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -89.59999999999999
$ws.Range("N4").Value = -530
$ws.Range("H22").Value = 197.8
$ws.Range("I22").Value = 197.8
$ws.Range("K22").Value = 197.8
$ws.Range("M22").Value = -24.80000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 947.7368
$ws.Range("I31").Value = 689.7941
$ws.Range("J31").Value = 1329.0435
$ws.Range("K31").Value = 689.7941
$ws.Range("L31").Value = 1329.0435
$ws.Range("M31").Value = -394.7941
$ws.Range("N31").Value = -1919.0435
$ws.Range("H34").Value = 947.7368
$ws.Range("I34").Value = 689.7941
$ws.Range("J34").Value = 1329.0435
$ws.Range("K34").Value = 689.7941
$ws.Range("L34").Value = 1329.0435
$ws.Range("M34").Value = -487.7941
$ws.Range("N34").Value = -1733.0435
$ws.Range("H58").Value = 901.6
$ws.Range("I58").Value = 941.2857
$ws.Range("K58").Value = 941.2857
$ws.Range("M58").Value = -738.2857
$ws.Range("H86").Value = 2788912.5
$ws.Range("J86").Value = 31559.285
$ws.Range("L86").Value = 31559.285
$ws.Range("N86").Value = -33805.285
$ws.Range("H89").Value = 2788912.5
$ws.Range("J89").Value = 31559.285
$ws.Range("L89").Value = 157796.425
$ws.Range("N89").Value = -169028.425
$ws.Range("H112").Value = 100000
$ws.Range("J112").Value = 100000
$ws.Range("L112").Value = 100000
$ws.Range("N112").Value = -102954
$ws.Range("H134").Value = 10754016
$ws.Range("I134").Value = 13334579
$ws.Range("J134").Value = 1669
$ws.Range("K134").Value = 40003737
$ws.Range("L134").Value = 5007
$ws.Range("M134").Value = -40001202
$ws.Range("N134").Value = -10077
$ws.Range("H136").Value = 901.6
$ws.Range("I136").Value = 941.2857
$ws.Range("K136").Value = 2823.8571
$ws.Range("M136").Value = -273.8571000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3573101
$ws.Range("J34").Value = 4168507.5
$ws.Range("L34").Value = 12505522.5
$ws.Range("N34").Value = -12505690.5
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2706
$ws.Range("H75").Value = 1160
$ws.Range("J75").Value = 1160
$ws.Range("L75").Value = 3480
$ws.Range("N75").Value = -5476
$ws.Range("H78").Value = 1160
$ws.Range("J78").Value = 1160
$ws.Range("L78").Value = 10440
$ws.Range("N78").Value = -20424
$ws.Range("H131").Value = 20834818
$ws.Range("I131").Value = 125000664
$ws.Range("J131").Value = 1650
$ws.Range("K131").Value = 375001992
$ws.Range("L131").Value = 4950
$ws.Range("M131").Value = -374996952
$ws.Range("N131").Value = -15030
$ws.Range("H140").Value = 29957.676
$ws.Range("I140").Value = 40226.793
$ws.Range("J140").Value = 2884.5454
$ws.Range("K140").Value = 120680.379
$ws.Range("L140").Value = 8653.636200000001
$ws.Range("M140").Value = -115500.379
$ws.Range("N140").Value = -19013.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1374.75
$ws.Range("I31").Value = 1374.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1374.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1082.75
$ws.Range("N31").Value = ""
$ws.Range("H37").Value = 1374.75
$ws.Range("I37").Value = 1374.75
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1374.75
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1097.75
$ws.Range("N37").Value = ""
$ws.Range("H126").Value = 3032.6667
$ws.Range("I126").Value = 2022.8572
$ws.Range("J126").Value = 3448.4707
$ws.Range("K126").Value = 6068.571599999999
$ws.Range("L126").Value = 10345.4121
$ws.Range("M126").Value = -3598.571599999999
$ws.Range("N126").Value = -15285.4121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2519.2
$ws.Range("I40").Value = 2519.2
$ws.Range("K40").Value = 2519.2
$ws.Range("M40").Value = -2383.2
$ws.Range("H46").Value = 7076.923
$ws.Range("I46").Value = 733.3333
$ws.Range("J46").Value = 8980
$ws.Range("K46").Value = 733.3333
$ws.Range("L46").Value = 8980
$ws.Range("M46").Value = -545.3333
$ws.Range("N46").Value = -9356
$ws.Range("H68").Value = 1527.5333
$ws.Range("I68").Value = 1314
$ws.Range("K68").Value = 1314
$ws.Range("M68").Value = -565
$ws.Range("H71").Value = 1527.5333
$ws.Range("I71").Value = 1314
$ws.Range("K71").Value = 6570
$ws.Range("M71").Value = -2826
$ws.Range("H93").Value = 2000
$ws.Range("J93").Value = 2000
$ws.Range("L93").Value = 2000
$ws.Range("N93").Value = -4496
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H136").Value = 9291
$ws.Range("I136").Value = 17972.334
$ws.Range("J136").Value = 1849.8572
$ws.Range("K136").Value = 53917.00199999999
$ws.Range("L136").Value = 5549.571599999999
$ws.Range("M136").Value = -51367.00199999999
$ws.Range("N136").Value = -10649.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40429
$ws.Range("J46").Value = 40429
$ws.Range("L46").Value = 40429
$ws.Range("N46").Value = -40891
$ws.Range("H100").Value = 894.6
$ws.Range("I100").Value = 618.25
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1236.5
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -695.5
$ws.Range("N100").Value = -5082
$ws.Range("H134").Value = 40429
$ws.Range("J134").Value = 40429
$ws.Range("L134").Value = 121287
$ws.Range("N134").Value = -126357
